$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: "prefab" field, loaded from the characters table
$ws.Range("G1").Value = "prefab"
$ws.Range("G2").Value = "string"
$ws.Range("G3").Value = "prefab_asset_key"
$ws.Range("G4").Value = "Level:Characters:ActorDragon"
$ws.Range("G5").Value = "Level:Characters:ActorDragon"

# Column width to match the other data columns (closest value the
# engine's character-width -> pixel-width rounding can reach to 30.125)
$ws.Range("G1").ColumnWidth = 29.375

# Header / type-row / display-row styling to match the existing table look
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").HorizontalAlignment = -4108

# Stray formatted (empty) cells left from selecting/formatting a larger range
$ws.Range("H5").HorizontalAlignment = -4108
$ws.Range("G13").HorizontalAlignment = -4108
